{"js": "// Update the date line at the top of the document.\nconst dateResults = context.document.body.search(\"2025-02-06 Thursday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"2025-02-07 Friday\", Word.InsertLocation.replace);\n}\n\n// Update the division problems in the table. Cells are addressed directly by\n// (row, column) so the two duplicate \"12\u00f75=\" cells each receive their own,\n// distinct replacement value, and existing run formatting (font/size) on\n// each cell is preserved.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst tbl = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, text: \"78\u00f74=\" }, // was 56\u00f77=\n  { row: 0, col: 1, text: \"52\u00f74=\" }, // was 49\u00f78=\n  { row: 0, col: 2, text: \"14\u00f78=\" }, // was 12\u00f75=\n  { row: 0, col: 3, text: \"83\u00f76=\" }, // was 69\u00f74=\n  { row: 0, col: 4, text: \"71\u00f76=\" }, // was 15\u00f73=\n\n  { row: 4, col: 0, text: \"19\u00f76=\" }, // was 68\u00f77=\n  { row: 4, col: 1, text: \"39\u00f72=\" }, // was 37\u00f73=\n  { row: 4, col: 2, text: \"55\u00f78=\" }, // was 91\u00f75=\n  { row: 4, col: 3, text: \"81\u00f79=\" }, // was 12\u00f75=\n  { row: 4, col: 4, text: \"22\u00f79=\" }, // was 63\u00f75=\n\n  { row: 8, col: 0, text: \"45\u00f77=\" }, // was 20\u00f72=\n  { row: 8, col: 1, text: \"69\u00f79=\" }, // was 89\u00f76=\n  { row: 8, col: 2, text: \"89\u00f78=\" }, // was 92\u00f73=\n  { row: 8, col: 3, text: \"50\u00f76=\" }, // was 90\u00f73=\n  { row: 8, col: 4, text: \"21\u00f78=\" }, // was 41\u00f74=\n\n  { row: 12, col: 0, text: \"84\u00f75=\" }, // was 25\u00f72=\n  { row: 12, col: 1, text: \"28\u00f78=\" }, // was 30\u00f72=\n  { row: 12, col: 2, text: \"47\u00f73=\" }, // was 37\u00f75=\n  { row: 12, col: 3, text: \"20\u00f78=\" }, // was 74\u00f75=\n  { row: 12, col: 4, text: \"27\u00f76=\" }, // was 43\u00f74=\n\n  { row: 16, col: 0, text: \"85\u00f76=\" }, // was 67\u00f72=\n  { row: 16, col: 1, text: \"15\u00f78=\" }, // was 88\u00f73=\n  { row: 16, col: 2, text: \"63\u00f73=\" }, // was 50\u00f72=\n  { row: 16, col: 3, text: \"38\u00f77=\" }, // was 93\u00f72=\n  { row: 16, col: 4, text: \"61\u00f78=\" }, // was 35\u00f73=\n];\n\nfor (const u of updates) {\n  const cell = tbl.getCell(u.row, u.col);\n  const range = cell.body.getRange();\n  range.insertText(u.text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line at the top of the document.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"2025-02-06 Thursday\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2025-02-07 Friday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Update the division problems in the table. Cells are addressed directly by\n# (row, column) so the two duplicate \"12\u00f75=\" cells get their own, distinct\n# replacement values.\n$tbl = $d.Tables.Item(1)\n\n$updates = @(\n    @{Row=1;  Col=1; Text=\"78\u00f74=\"},\n    @{Row=1;  Col=2; Text=\"52\u00f74=\"},\n    @{Row=1;  Col=3; Text=\"14\u00f78=\"},\n    @{Row=1;  Col=4; Text=\"83\u00f76=\"},\n    @{Row=1;  Col=5; Text=\"71\u00f76=\"},\n\n    @{Row=5;  Col=1; Text=\"19\u00f76=\"},\n    @{Row=5;  Col=2; Text=\"39\u00f72=\"},\n    @{Row=5;  Col=3; Text=\"55\u00f78=\"},\n    @{Row=5;  Col=4; Text=\"81\u00f79=\"},\n    @{Row=5;  Col=5; Text=\"22\u00f79=\"},\n\n    @{Row=9;  Col=1; Text=\"45\u00f77=\"},\n    @{Row=9;  Col=2; Text=\"69\u00f79=\"},\n    @{Row=9;  Col=3; Text=\"89\u00f78=\"},\n    @{Row=9;  Col=4; Text=\"50\u00f76=\"},\n    @{Row=9;  Col=5; Text=\"21\u00f78=\"},\n\n    @{Row=13; Col=1; Text=\"84\u00f75=\"},\n    @{Row=13; Col=2; Text=\"28\u00f78=\"},\n    @{Row=13; Col=3; Text=\"47\u00f73=\"},\n    @{Row=13; Col=4; Text=\"20\u00f78=\"},\n    @{Row=13; Col=5; Text=\"27\u00f76=\"},\n\n    @{Row=17; Col=1; Text=\"85\u00f76=\"},\n    @{Row=17; Col=2; Text=\"15\u00f78=\"},\n    @{Row=17; Col=3; Text=\"63\u00f73=\"},\n    @{Row=17; Col=4; Text=\"38\u00f77=\"},\n    @{Row=17; Col=5; Text=\"61\u00f78=\"}\n)\n\nforeach ($u in $updates) {\n    $cell = $tbl.Cell($u.Row, $u.Col)\n    $cellRange = $cell.Range\n    $cellRange.MoveEnd(1, -1) | Out-Null\n    $cellRange.Text = $u.Text\n}\n"}
